$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country rankings / stats (re-sorted data refresh) ---
$ws.Cells.Item(30, 1).Value = "Chequia"
$ws.Cells.Item(30, 2).Value = 5589
$ws.Cells.Item(30, 3).Value = 20
$ws.Cells.Item(30, 4).Value = 309
$ws.Cells.Item(30, 5).Value = 5167
$ws.Cells.Item(30, 6).Value = 96
$ws.Cells.Item(30, 7).Value = 1
$ws.Cells.Item(30, 8).Value = 113

$ws.Cells.Item(31, 1).Value = "Polonia"
$ws.Cells.Item(31, 2).Value = 5575
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 284
$ws.Cells.Item(31, 5).Value = 5117
$ws.Cells.Item(31, 6).Value = 160
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 174

$ws.Cells.Item(51, 1).Value = "Ucrania"
$ws.Cells.Item(51, 2).Value = 2203
$ws.Cells.Item(51, 3).Value = 311
$ws.Cells.Item(51, 4).Value = 61
$ws.Cells.Item(51, 5).Value = 2073
$ws.Cells.Item(51, 6).Value = 33
$ws.Cells.Item(51, 7).Value = 12
$ws.Cells.Item(51, 8).Value = 69

$ws.Cells.Item(52, 1).Value = "Grecia"
$ws.Cells.Item(52, 2).Value = 1955
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 269
$ws.Cells.Item(52, 5).Value = 1599
$ws.Cells.Item(52, 6).Value = 79
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 87

$ws.Cells.Item(53, 1).Value = "Sudafrica"
$ws.Cells.Item(53, 2).Value = 1934
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 95
$ws.Cells.Item(53, 5).Value = 1821
$ws.Cells.Item(53, 6).Value = 7
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 18

$ws.Cells.Item(54, 1).Value = "Singapur"
$ws.Cells.Item(54, 2).Value = 1910
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 460
$ws.Cells.Item(54, 5).Value = 1444
$ws.Cells.Item(54, 6).Value = 29
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 6

$ws.Cells.Item(55, 1).Value = "Argentina"
$ws.Cells.Item(55, 2).Value = 1894
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(55, 4).Value = 365
$ws.Cells.Item(55, 5).Value = 1450
$ws.Cells.Item(55, 6).Value = 96
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 79

$ws.Cells.Item(90, 1).Value = "Oman"
$ws.Cells.Item(90, 2).Value = 484
$ws.Cells.Item(90, 3).Value = 27
$ws.Cells.Item(90, 4).Value = 109
$ws.Cells.Item(90, 5).Value = 372
$ws.Cells.Item(90, 6).Value = 3
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 3

$ws.Cells.Item(91, 1).Value = "Uruguay"
$ws.Cells.Item(91, 2).Value = 473
$ws.Cells.Item(91, 3).Value = 17
$ws.Cells.Item(91, 4).Value = 206
$ws.Cells.Item(91, 5).Value = 260
$ws.Cells.Item(91, 6).Value = 13
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 7

$ws.Cells.Item(112, 1).Value = "Georgia"
$ws.Cells.Item(112, 2).Value = 227
$ws.Cells.Item(112, 3).Value = 9
$ws.Cells.Item(112, 4).Value = 52
$ws.Cells.Item(112, 5).Value = 172
$ws.Cells.Item(112, 6).Value = 6
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 3

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 08:52"
